$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: 2025-03-17, Amna (member 4), 01:47:07 - 01:47:16
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "2025-03-17"
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "Amna"
$ws.Range("D13").Value = "01:47:07"
$ws.Range("E13").Value = "01:47:16"

# Row 14: 2025-03-17, nabeel (member 2), 01:47:13 - 01:47:20
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2025-03-17"
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "nabeel"
$ws.Range("D14").Value = "01:47:13"
$ws.Range("E14").Value = "01:47:20"
